$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (after the weekly shuffle), columns: D, M, N, O, P, R, S
$data = @{
    2  = @{ D = 44238; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 }
    3  = @{ D = 44174; M = 200; N = 3200; O = 3200; P = 3200; R = "Provincia de Curicó"; S = 1600 }
    4  = @{ D = 44188; M = 150; N = 3000; O = 3400; P = 3240; R = "Provincia de Linares"; S = 1620 }
    5  = @{ D = 44617; M = 90;  N = 6500; O = 6500; P = 6500; R = "Provincia de Curicó"; S = 3250 }
    7  = @{ D = 44236; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 }
    8  = @{ D = 44208; M = 85;  N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 }
    9  = @{ D = 44586; M = 250; N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó"; S = 2500 }
    10 = @{ D = 44237; M = 100; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó"; S = 1900 }
    11 = @{ D = 44231; M = 150; N = 3400; O = 3400; P = 3400; R = "Provincia de Curicó"; S = 1700 }
    12 = @{ D = 44168; M = 170; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    14 = @{ D = 44582; M = 380; N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó"; S = 2500 }
    15 = @{ D = 44533; M = 150; N = 4000; O = 4000; P = 4000; R = "Provincia de Curicó"; S = 2000 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 13).Value = $vals.M
    $ws.Cells.Item($row, 14).Value = $vals.N
    $ws.Cells.Item($row, 15).Value = $vals.O
    $ws.Cells.Item($row, 16).Value = $vals.P
    $ws.Cells.Item($row, 18).Value = $vals.R
    $ws.Cells.Item($row, 19).Value = $vals.S
}
